# Update Contact_List: swap in a new set of contacts, add mailto hyperlinks
# on the Email column, and refresh the selection / column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ellipsis = [char]0x2026

# --- text edits, issued in the same order the original author typed them ---
$ws.Cells.Item(2, 3).Value = "Help, I've fallen and I cannot get up"
$ws.Cells.Item(3, 3).Value = "How can you have any pudding, if you don't eat your meat?"
$ws.Cells.Item(3, 1).Value = "Roger Gilmour"
$ws.Cells.Item(2, 1).Value = "Sam Power"
$ws.Cells.Item(4, 3).Value = "Thanks for the fast delivery, you guys rock!!!"
$ws.Cells.Item(9, 3).Value = "Seriously$($ellipsis)????"

# --- email column: replace address + attach a mailto hyperlink ---
$ws.Hyperlinks.Add($ws.Cells.Item(2, 2), "mailto:SamPower@Null.Org", "", "", "SamPower@Null.Org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 2), "mailto:RogerGilmour@Null.Org", "", "", "RogerGilmour@Null.Org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 2), "mailto:DanLagomarsino@Null.Org", "", "", "DanLagomarsino@Null.Org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 2), "mailto:GilStrong@Null.org", "", "", "GilStrong@Null.org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 2), "mailto:KasperNash@Null.org", "", "", "KasperNash@Null.org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 2), "mailto:CainCombs@Null.org", "", "", "CainCombs@Null.org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 2), "mailto:WingCollier@Null.org", "", "", "WingCollier@Null.org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 2), "mailto:DianaShaw@Null.org", "", "", "DianaShaw@Null.org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 2), "mailto:ElijahBarnett@Null.org", "", "", "ElijahBarnett@Null.org") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 2), "mailto:HermanCrosby@Null.org", "", "", "HermanCrosby@Null.org") | Out-Null

# --- column widths follow the new, wider content ---
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# --- selection moves to B12, matching the saved view state ---
$ws.Range("B12").Select()
